$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces, followed by three red (C00000) runs that read
#    "(This is a change – Version for branch alternate)"
# ---------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$insertPos = $p1.End - 1

$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter("  ")

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r2.Font.Color = 192

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter("rsion for branch alternate")
$r3.Font.Color = 192

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(")")
$r4.Font.Color = 192

# ---------------------------------------------------------------
# 2) "Crispian's Day speech from Shakespear's Henry V [Source -
#    Wikipedia]" paragraph: merge run boundaries / move proofErr
#    markers without altering the visible text.
# ---------------------------------------------------------------
$p4 = $d.Paragraphs(4).Range
$base = $p4.Start

# "...Day speech from" + " " -> "...Day speech from "
$rSpace = $d.Range($base + 26, $base + 27)
$rSpace.Delete()
$rFrom = $d.Range($base + 26, $base + 26)
$rFrom.InsertAfter(" ")

# " Henry V" + " " + "[Source " + "-" + " Wikipedia" + "]"
#   -> " Henry V [Source - Wikipedia]"
$rHenry = $d.Range($base + 47, $base + 47)
$rHenry.InsertAfter(" [Source " + [char]0x2013 + " Wikipedia]")
$tailStart = $rHenry.End
$rTail = $d.Range($tailStart, $tailStart + 21)
$rTail.Delete()

# ---------------------------------------------------------------
# 3) Add two new paragraphs at the very end of the document body.
# ---------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$paraCount = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs($paraCount - 1).Range
$secondToLast.Style = "larger"
$secondToLast.ParagraphFormat.Shading.BackgroundPatternColor = 16777215
$secondToLast.ParagraphFormat.SpaceBefore = 0
$secondToLast.ParagraphFormat.SpaceAfter = 7.5
